# Update scripts with new TPM values (Efnb3-Ephb2, YoungD0)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "ECs" sending-cluster rows (original rows 2 and 3).
# The remaining rows (originally 4-7) shift up to become rows 2-5, and
# their Sending/Target cluster labels already match the desired output.
$ws.Rows("2:3").Delete()

# Row 2 (FAPs -> FAPs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1498043333333333
$ws.Range("H2").Value = 0.449413
$ws.Range("I2").Value = 0.08722868471333377
$ws.Range("J2").Value = 0.08722868471333377
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.346253666666667
$ws.Range("N2").Value = 19.038761
$ws.Range("O2").Value = 0.921725411846598
$ws.Range("P2").Value = 0.9217254118465981
$ws.Range("Q2").Value = 0.9506962996992223
$ws.Range("R2").Value = 8.556266697293001
$ws.Range("S2").Value = 0.08040089534223462
$ws.Range("T2").Value = 0.08040089534223462

# Row 3 (FAPs -> MuSCs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1498043333333333
$ws.Range("H3").Value = 0.449413
$ws.Range("I3").Value = 0.08722868471333377
$ws.Range("J3").Value = 0.08722868471333377
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5389353333333333
$ws.Range("N3").Value = 1.616806
$ws.Range("O3").Value = 0.07827458815340194
$ws.Range("P3").Value = 0.07827458815340194
$ws.Range("Q3").Value = 0.08073484831977779
$ws.Range("R3").Value = 0.726613634878
$ws.Range("S3").Value = 0.006827789371099148
$ws.Range("T3").Value = 0.006827789371099148

# Row 4 (MuSCs -> FAPs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.567570333333333
$ws.Range("H4").Value = 4.702711
$ws.Range("I4").Value = 0.9127713152866662
$ws.Range("J4").Value = 0.9127713152866662
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.346253666666667
$ws.Range("N4").Value = 19.038761
$ws.Range("O4").Value = 0.921725411846598
$ws.Range("P4").Value = 0.9217254118465981
$ws.Range("Q4").Value = 9.948198975674556
$ws.Range("R4").Value = 89.533790781071
$ws.Range("S4").Value = 0.8413245165043634
$ws.Range("T4").Value = 0.8413245165043635

# Row 5 (MuSCs -> MuSCs)
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.567570333333333
$ws.Range("H5").Value = 4.702711
$ws.Range("I5").Value = 0.9127713152866662
$ws.Range("J5").Value = 0.9127713152866662
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5389353333333333
$ws.Range("N5").Value = 1.616806
$ws.Range("O5").Value = 0.07827458815340194
$ws.Range("P5").Value = 0.07827458815340194
$ws.Range("Q5").Value = 0.8448190401184444
$ws.Range("R5").Value = 7.603371361065999
$ws.Range("S5").Value = 0.07144679878230278
$ws.Range("T5").Value = 0.07144679878230278
